$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 50, shifting existing rows 50:75 down to 51:76
$ws.Rows.Item(50).Insert()

# Populate the newly inserted row 50 with the new weekly record
$ws.Cells.Item(50, 1).Value = 11
$ws.Cells.Item(50, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(50, 3).Value = "Bíobío"
$ws.Cells.Item(50, 4).Value = 44825
$ws.Cells.Item(50, 5).Value = 8
$ws.Cells.Item(50, 6).Value = 100112013
$ws.Cells.Item(50, 7).Value = "Alcachofa"
$ws.Cells.Item(50, 8).Value = "Argentina(o)"
$ws.Cells.Item(50, 9).Value = "Primera"
$ws.Cells.Item(50, 10).Value = 60
$ws.Cells.Item(50, 11).Value = 12000
$ws.Cells.Item(50, 12).Value = 13000
$ws.Cells.Item(50, 13).Value = 12500
$ws.Cells.Item(50, 14).Value = "$/caja 50 unidades"
$ws.Cells.Item(50, 15).Value = "Provincia de Limarí"
$ws.Cells.Item(50, 16).Value = 250
$ws.Cells.Item(50, 17).Value = 50
$ws.Cells.Item(50, 18).Value = "Hortaliza"
